$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.873.27"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "2.498.53"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.31"
$ws.Range("E5").Value = "  +2.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.72"
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("E8").Value = "  +2.39%  "
$ws.Range("D9").Value = "2.500.09"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0991"
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("E11").Value = "  -2.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.14"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.329"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").Value = "2.946.83"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "58.698.95"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.29"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "2.504.07"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.58"
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.25"
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.01"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.16"
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.79"
$ws.Range("E24").Value = "  +3.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.408"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.159"
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.42"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.63"
$ws.Range("E29").Value = "  +2.36%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0754"
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.26"
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.17"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.07"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.22"
$ws.Range("E37").Value = "  -4.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.95"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.51"
$ws.Range("E39").Value = "  +3.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.824"
$ws.Range("E40").Value = "  +6.01%  "
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "273.73"
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "131.13"
$ws.Range("E44").Value = "  +6.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.01"
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.592"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0932"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0509"
$ws.Range("E48").Value = "  +2.35%  "
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.70"
$ws.Range("E50").Value = "  -2.21%  "
$ws.Range("D51").Value = "1.750.27"
$ws.Range("E51").Value = "  +0.39%  "
